$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.574.60"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "1.923.23"
$ws.Range("E3").Value = "  +4.03%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.08"
$ws.Range("E5").Value = "  +4.88%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4720"
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2890"
$ws.Range("E8").Value = "  +4.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06778"
$ws.Range("E9").Value = "  +5.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "105.19"
$ws.Range("E10").Value = "  +8.24%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "1.916.18"
$ws.Range("E12").Value = "  +4.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07703"
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("E14").Value = "  +6.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6721"
$ws.Range("E15").Value = "  +7.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "290.15"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "30.584.13"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007606"
$ws.Range("E18").Value = "  +3.37%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.91"
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("D21").Value = "2.166.01"
$ws.Range("E21").Value = "  +4.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.468"
$ws.Range("E22").Value = "  +9.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.330"
$ws.Range("E24").Value = "  +4.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.417"
$ws.Range("E25").Value = "  +3.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.73"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.08"
$ws.Range("E27").Value = "  +9.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.117"
$ws.Range("E28").Value = "  +9.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1075"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.391"
$ws.Range("E30").Value = "  +4.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.183"
$ws.Range("E31").Value = "  +5.03%  "
$ws.Range("E32").Value = "  +8.27%  "
$ws.Range("E33").Value = "  +3.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7436"
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.153"
$ws.Range("E35").Value = "  +3.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02077"
$ws.Range("E36").Value = "  +10.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.747"
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.050"
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "111.57"
$ws.Range("E40").Value = "  +5.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8807"
$ws.Range("E41").Value = "  +3.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4395"
$ws.Range("E42").Value = "  +8.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.890"
$ws.Range("E43").Value = "  +3.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "67.08"
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.219"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.279"
$ws.Range("E47").Value = "  +3.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.10"
$ws.Range("E48").Value = "  +19.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1238"
$ws.Range("E49").Value = "  +4.56%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4093"
$ws.Range("E50").Value = "  +10.68%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.92"
$ws.Range("E51").Value = "  +3.35%  "
